# Updated sampling ranges for experiment and added climate change factor to hydropower.
#
# For every data row on the active sheet (rows 2-226), the "max_35" (column H)
# and "min_35" (column I) scaling-range bounds are reset to 1, collapsing the
# previous per-variable multiplicative ranges to a uniform factor of 1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 8).End(-4162).Row  # xlUp = -4162, column H

if ($lastRow -lt 2) {
    $lastRow = 226
}

$ws.Range("H2:I$lastRow").Value = 1
